# Updates the cryptocurrency price ("D") and 1h volume-change ("E")
# columns on Sheet1 to match the latest scraped snapshot.
# The Price column holds plain text (values like "26.784.46" use
# dots as thousands separators, not valid numbers) so a leading
# apostrophe forces Excel to store the literal text instead of
# auto-converting look-alike numbers (e.g. "0.503") to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.784.46"
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = "'1.644.49"
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("E4").Value = '  +0.66%  '
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("D6").Value = "'0.503"
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("E7").Value = '  +0.72%  '
$ws.Range("D8").Value = "'0.252"
$ws.Range("E8").Value = '  -0.85%  '
$ws.Range("E9").Value = '  -0.50%  '
$ws.Range("D10").Value = "'19.19"
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("D12").Value = "'1.869.17"
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = "'1.640.58"
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = "'4.18"
$ws.Range("E14").Value = '  -0.86%  '
$ws.Range("D15").Value = "'0.527"
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("D16").Value = "'64.60"
$ws.Range("E16").Value = '  -3.01%  '
$ws.Range("D17").Value = "'26.790.85"
$ws.Range("E18").Value = '  -2.47%  '
$ws.Range("D19").Value = "'214.49"
$ws.Range("E19").Value = '  -3.03%  '
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("E21").Value = '  -1.01%  '
$ws.Range("D22").Value = "'2.40"
$ws.Range("E22").Value = '  +11.95%  '
$ws.Range("E23").Value = '  -1.47%  '
$ws.Range("D24").Value = "'9.36"
$ws.Range("E24").Value = '  -2.32%  '
$ws.Range("D25").Value = "'144.83"
$ws.Range("E25").Value = '  -1.51%  '
$ws.Range("E26").Value = '  +1.04%  '
$ws.Range("D27").Value = "'0.119"
$ws.Range("E27").Value = '  -2.34%  '
$ws.Range("D28").Value = "'7.12"
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").Value = "'15.69"
$ws.Range("E29").Value = '  -1.61%  '
$ws.Range("D30").Value = "'0.0514"
$ws.Range("E30").Value = '  -1.75%  '
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("D32").Value = "'3.33"
$ws.Range("E32").Value = '  -3.22%  '
$ws.Range("E33").Value = '  -2.65%  '
$ws.Range("D34").Value = "'1.287.41"
$ws.Range("E34").Value = '  -0.37%  '
$ws.Range("D35").Value = "'1.54"
$ws.Range("E35").Value = '  -2.27%  '
$ws.Range("D36").Value = "'2.43"
$ws.Range("E36").Value = '  +1.38%  '
$ws.Range("E37").Value = '  -4.37%  '
$ws.Range("D38").Value = "'0.541"
$ws.Range("E38").Value = '  +2.33%  '
$ws.Range("D39").Value = "'0.826"
$ws.Range("E39").Value = '  -0.82%  '
$ws.Range("E40").Value = '  +0.70%  '
$ws.Range("D41").Value = "'0.813"
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").Value = "'5.35"
$ws.Range("D44").Value = "'1.794.92"
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").Value = "'91.46"
$ws.Range("E45").Value = '  -2.51%  '
$ws.Range("D46").Value = "'60.11"
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("E47").Value = '  -1.00%  '
$ws.Range("E48").Value = '  -1.23%  '
$ws.Range("D49").Value = "'0.0519"
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("E50").Value = '  -2.09%  '
$ws.Range("E51").Value = '  -0.52%  '
